# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps on the Overview, zh-cn and de-de
# sheets to reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 36ab7ac0-3d59-48ce-873e-566bba4bdaa3.md row
$wsOverview.Range("G4").Value = "2016-08-27 02:44:35"

# zh-cn sheet: Correspond Handoff / Handback datetimes for the same row
$wsZhCn.Range("H4").Value = "2016-08-27 02:44:31"
$wsZhCn.Range("K4").Value = "2016-08-27 02:44:55"

# de-de sheet: Correspond Handoff / Handback datetimes for the same row
$wsDeDe.Range("H4").Value = "2016-08-27 02:44:35"
$wsDeDe.Range("K4").Value = "2016-08-27 02:45:06"
